{"js": "// \"no more FE split, add l4\"\n// Change the Nuro job-title line from:\n//   Software Engineer<TAB>Mountain View, CA\n// to:\n//   Software Engineer L3 \u2192 L4<TAB>Mountain View, CA\n//\n// We locate the unique occurrence of \"Software Engineer\" that is the title\n// on the Nuro line (i.e. immediately followed by a tab, to avoid matching\n// \"Software Engineer Intern\" used at Braze/other jobs), then insert the new\n// \" L3 \", \"\u2192\", \" L4\" text right after it.\n\nconst body = context.document.body;\nconst results = body.search(\"Software Engineer\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\n// Load the first paragraph for every hit so we can disambiguate which one is\n// the Nuro job-title line (text starts with \"Software Engineer\" + TAB).\nfor (let i = 0; i < results.items.length; i++) {\n  const par = results.items[i].paragraphs.getFirst();\n  par.load(\"text\");\n  results.items[i]._titlePara = par;\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < results.items.length; i++) {\n  if (results.items[i]._titlePara.text.indexOf(\"Software Engineer\\t\") === 0) {\n    target = results.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Software Engineer' (Nuro) title run\");\n}\n\n// Insert the three new text pieces right after \"Software Engineer\".\nlet afterRange = target.insertText(\" L3 \", \"After\");\nawait context.sync();\nafterRange = afterRange.insertText(\"\\u2192\", \"After\");\nawait context.sync();\nafterRange.insertText(\" L4\", \"After\");\nawait context.sync();\n", "ps1": "# \"no more FE split, add l4\"\n# Change the Nuro job-title line from:\n#   Software Engineer<TAB>Mountain View, CA\n# to:\n#   Software Engineer L3 -> L4<TAB>Mountain View, CA\n#\n# Locate the paragraph whose text is exactly \"Software Engineer<TAB>Mountain\n# View, CA\" (the Nuro job-title line; other \"Software Engineer Intern\" lines\n# at Braze etc. won't match), then insert the new \" L3 \", arrow, \" L4\" text\n# right after \"Software Engineer\" there.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"Software Engineer`tMountain View, CA`r\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $fr = $target.Range.Duplicate\n    $fr.Find.ClearFormatting()\n    $fr.Find.MatchCase = $true\n    $fr.Find.Execute(\"Software Engineer\") | Out-Null\n\n    if ($fr.Find.Found) {\n        $ins = $fr.Duplicate\n        $ins.Collapse(0)              # wdCollapseEnd\n        $ins.InsertAfter(\" L3 \")\n\n        $fr2 = $target.Range.Duplicate\n        $fr2.Find.ClearFormatting()\n        $fr2.Find.MatchCase = $true\n        $fr2.Find.Execute(\" L3 \") | Out-Null\n        $ins2 = $fr2.Duplicate\n        $ins2.Collapse(0)\n        $ins2.InsertAfter([char]0x2192)\n\n        $fr3 = $target.Range.Duplicate\n        $fr3.Find.ClearFormatting()\n        $fr3.Find.MatchCase = $true\n        $fr3.Find.Execute([char]0x2192) | Out-Null\n        $ins3 = $fr3.Duplicate\n        $ins3.Collapse(0)\n        $ins3.InsertAfter(\" L4\")\n    }\n}\n"}
